# Apply "adding my version of proposal v1":
#  1. Remove the older brainstormed ideas (everything from "Online checker for
#     change in privacy laws." through the blank paragraph that precedes the
#     "User gets asked ..." proposal), leaving only the "User gets asked ..."
#     write-up in place.
#  2. Tweak the wording of the remaining proposal: "... checklist is
#     satisfied and it does highlight warnings beforehand." becomes
#     "... checklist is satisfied, and it does highlight warnings
#     beforehand." (comma added after "satisfied").

$d = $word.ActiveDocument

# --- Step 1: delete the obsolete idea paragraphs -----------------------
$startRange = $d.Content
$startRange.Find.Execute("Online checker for change in privacy laws.")

$endRange = $d.Content
$endRange.Find.Execute("User gets asked what type of data")

$obsolete = $d.Range($startRange.Start, $endRange.Start)
$obsolete.Delete()

# --- Step 2: insert the missing comma -----------------------------------
$satisfied = $d.Content
$satisfied.Find.Execute("satisfied")
$insertionPoint = $d.Range($satisfied.End, $satisfied.End)
$insertionPoint.InsertAfter(",")
